$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 11 ("Closing Comments & Asks") - Content Placeholder 2
# ---------------------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$shp11 = $s11.Shapes.Item(2)
$tf11 = $shp11.TextFrame
$tr11 = $tf11.TextRange

# Remove the "TLS is not just for the Web." and "A compression certificate
# dictionary is used in draft-ietf-tls-ctls as well." paragraphs entirely
# (including their trailing paragraph mark).
$full = $tr11.Text
$start = $full.IndexOf("TLS is not just for the Web.")
$endMarker = "A compression certificate dictionary is used in draft-ietf-tls-ctls as well."
$endIdx = $full.IndexOf($endMarker)
$delLen = ($endIdx + $endMarker.Length) - $start + 1
$delRange = $tr11.Characters($start + 1, $delLen)
$delRange.Delete()

# "But also " -> "But also, "
$full = $tr11.Text
$idx = $full.IndexOf("But also ")
$sub = $tr11.Characters($idx + 1, "But also ".Length)
$sub.Text = "But also, "

# "let's not forget" -> "let’s not forget, TLS is not just for the Web."
$full = $tr11.Text
$searchTarget = "let's not forget"
$idx = $full.IndexOf($searchTarget)
$sub = $tr11.Characters($idx + 1, $searchTarget.Length)
$sub.Text = "let’s not forget, TLS is not just for the Web."

# Recompute autofit now that there is less text (drops the stale
# lnSpcReduction="10000" -> normAutofit with no reduction).
$tf11.AutoSize = 2

# ---------------------------------------------------------------------------
# Slide 9 ("About ICA lists") - Content Placeholder 2
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$shp9 = $s9.Shapes.Item(2)
$tf9 = $shp9.TextFrame
$tr9 = $tf9.TextRange

# Start a new paragraph (inherits lvl="1" from the preceding paragraph) and
# add the new sentence as several runs, then a trailing empty paragraph.
[void]$tr9.InsertAfter("`rdraft-")
[void]$tr9.InsertAfter("ietf")
[void]$tr9.InsertAfter("-")
[void]$tr9.InsertAfter("tls-ctls")
[void]$tr9.InsertAfter(" also uses a compression certificate ")
[void]$tr9.InsertAfter("dictionar")
[void]$tr9.InsertAfter(". `r")
